$d = $word.ActiveDocument

# 1. Insert two new paragraphs right before the paragraph that contains
#    "192 horas /16 semanas = 12 horas a la semana"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "192 horas /16 semanas = 12 horas a la semana"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "1 credito equivale ha 48 horas  el semestre ^p48 horas /16 semanas  = 3 horas^p192 horas /16 semanas = 12 horas a la semana"
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
